$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New sensor log rows to append below the existing A1:F12 data.
$newRows = @(
    @("2026-02-01", "17:13:23", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:13:32", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:13:43", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:14:14", "17:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-02-01", "17:14:25", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:14:36", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:14:46", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:14:56", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 13
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a plain "yyyy-mm-dd" text label (matches the existing
    # rows above it), so force Text formatting before assigning the value
    # to stop Excel from auto-converting it into a date serial number.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
